$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: delete B2 and D2 entirely, update C2 and E2
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -7.3733077377562868
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -8.2961493261296724

# Row 3: value tweaks
$ws.Range("B3").Value = -10.616310651571711
$ws.Range("C3").Value = -1.8565732042816296
$ws.Range("D3").Value = -15.401338659918059
$ws.Range("E3").Value = 20.513778672349321

# Update the selection to match the post-edit state
$ws.Range("B1:E3").Select()
